$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dimension/content: rows 6-21, columns A-C (Job, Machine, Duration)
$data = @(
    @(1,2,3),
    @(1,3,7),
    @(2,1,3),
    @(2,3,3),
    @(2,3,1),
    @(2,2,8),
    @(2,1,3),
    @(2,3,1),
    @(3,2,2),
    @(3,1,4),
    @(3,3,1),
    @(3,1,2),
    @(4,4,2),
    @(4,4,6),
    @(4,3,2),
    @(4,3,1)
)

$startRow = 6
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
}
